$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.845.79'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.640.97'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.90'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.79'
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.867.54'
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.640.14'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.24'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.874.08'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.48'
$ws.Range('E20').Value = '  +2.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.33'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.00'
$ws.Range('E22').Value = '  +0.42%  '
$ws.Range('E23').Value = '  +2.57%  '
$ws.Range('E24').Value = '  +3.91%  '
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.77'
$ws.Range('E26').Value = '  +3.11%  '
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.98'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.59'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0497'
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.911'
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.134.67'
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.547'
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.56'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.73'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.810'
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.776.50'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0108'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.54'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('E49').Value = '  +5.90%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  -1.96%  '
